$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 19, shifting existing rows 19-30 down to 21-32.
$ws.Range("A19:T20").EntireRow.Insert()

# --- New row 19 ---
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C19").Value = "Arica y Parinacota"
$ws.Range("D19").Value = 45264
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100107
$ws.Range("H19").Value = "Otros"
$ws.Range("I19").Value = 100107002
$ws.Range("J19").Value = "Chirimoya"
$ws.Range("K19").Value = "Cultivar IV Región"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 90
$ws.Range("N19").Value = 20000
$ws.Range("O19").Value = 20000
$ws.Range("P19").Value = 20000
$ws.Range("Q19").Value = "$/bandeja 10 kilos"
$ws.Range("R19").Value = "Región de Coquimbo"
$ws.Range("S19").Value = 2000
$ws.Range("T19").Value = 10

# --- New row 20 ---
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C20").Value = "Arica y Parinacota"
$ws.Range("D20").Value = 45264
$ws.Range("E20").Value = 15
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100107
$ws.Range("H20").Value = "Otros"
$ws.Range("I20").Value = 100107002
$ws.Range("J20").Value = "Chirimoya"
$ws.Range("K20").Value = "Cultivar IV Región"
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 22000
$ws.Range("O20").Value = 22000
$ws.Range("P20").Value = 22000
$ws.Range("Q20").Value = "$/bandeja 10 kilos"
$ws.Range("R20").Value = "Región de Coquimbo"
$ws.Range("S20").Value = 2200
$ws.Range("T20").Value = 10
